$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 41 and 42: coin identity swap (Stellar <-> EnergySwap) ---
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.06%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.113"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.66%  "

# --- Remaining price / volume refreshes ---
$ws.Range("D2").Value = "50.043.61"
$ws.Range("E2").Value = "  +4.39%  "

$ws.Range("D3").Value = "2.676.46"
$ws.Range("E3").Value = "  +7.81%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "326.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.530"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.23%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.558"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.19%  "

$ws.Range("E10").Value = "  +5.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.09%  "

$ws.Range("E12").Value = "  +3.46%  "

$ws.Range("E13").Value = "  +0.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.14%  "

$ws.Range("D15").Value = "3.090.97"
$ws.Range("E15").Value = "  +7.72%  "

$ws.Range("D16").Value = "2.701.56"
$ws.Range("E16").Value = "  +8.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.878"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.75%  "

$ws.Range("D18").Value = "49.991.20"
$ws.Range("E18").Value = "  +4.48%  "

$ws.Range("E19").Value = "  +4.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.74%  "

$ws.Range("E22").Value = "  +3.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "278.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("E25").Value = "  +4.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.71%  "

$ws.Range("E29").Value = "  -1.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.23%  "

$ws.Range("E31").Value = "  +4.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.32%  "

$ws.Range("E33").Value = "  +4.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0820"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.30%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("E38").Value = "  +7.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.21%  "

$ws.Range("E43").Value = "  +0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0320"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.85%  "

$ws.Range("D45").Value = "2.119.50"
$ws.Range("E45").Value = "  +6.28%  "

$ws.Range("E46").Value = "  +5.83%  "

$ws.Range("E47").Value = "  +14.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.72%  "

$ws.Range("E50").Value = "  +5.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.85"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.88%  "

